$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.117350816726685
$ws.Range("B1").Value = 2.285988330841064
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.493206262588501
$ws.Range("E1").Value = 0.97004234790802
